$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column F (dSF) values per diff
$ws.Range("F5").Value = -2
$ws.Range("F6").Value = 2
$ws.Range("F10").Value = 8
$ws.Range("F11").Value = -7
$ws.Range("F12").Value = -2
